$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# Header text updates: Volume Number 37 -> 38,
# Week covering 9/11/2023-9/17/2023 -> 9/18/2023-9/24/2023.
# Use Characters() so only the specific run text changes (keeps
# the rest of the rich-text run formatting intact).
# ============================================================
$volRange = $ws.Range("A8")
$volRange.Characters(21, 2).Text = "38"

$weekRange = $ws.Range("C9")
$weekRange.Characters(27, 9).Text = "9/18/2023"
$weekRange.Characters(47, 9).Text = "9/24/2023"

# ============================================================
# Crime-stat table (rows 14-30): refreshed weekly figures.
# A few cells flip between a numeric 0 and the literal "0"/
# "***.*" placeholder text used elsewhere in the sheet for a
# blank/undefined percentage. Copy a same-column donor cell
# that already has the exact target style + content so the
# cell's type (t="s" vs numeric) and style id match exactly,
# instead of letting Excel invent a brand-new style.
# ============================================================
$ws.Range("C23").Copy($ws.Range("C18"))
$ws.Range("F28").Copy($ws.Range("G28"))
$ws.Range("E28").Copy($ws.Range("H28"))
$ws.Range("F29").Copy($ws.Range("G29"))
$ws.Range("E29").Copy($ws.Range("H29"))
$ws.Range("F30").Copy($ws.Range("D30"))
$ws.Range("K30").Copy($ws.Range("E30"))
$ws.Range("F30").Copy($ws.Range("G30"))
$ws.Range("K30").Copy($ws.Range("H30"))

# ============================================================
# Final values for every changed cell in the table.
# ============================================================
$ws.Range("N14").Value = -86.666666666666
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 29
$ws.Range("K15").Value = -44.827586206896
$ws.Range("L15").Value = -15.78947368421
$ws.Range("N15").Value = -30.434782608695
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = -17.142857142857
$ws.Range("I16").Value = 243
$ws.Range("J16").Value = 215
$ws.Range("K16").Value = 13.023255813953
$ws.Range("L16").Value = 38.068181818181
$ws.Range("M16").Value = -3.95256916996
$ws.Range("N16").Value = -75.204081632653
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -11.111111111111
$ws.Range("F17").Value = 42
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 354
$ws.Range("J17").Value = 317
$ws.Range("K17").Value = 11.67192429022
$ws.Range("L17").Value = 27.797833935018
$ws.Range("M17").Value = 41.6
$ws.Range("N17").Value = 9.597523219814
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -38.461538461538
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = -0.9009009009
$ws.Range("L18").Value = 22.222222222222
$ws.Range("M18").Value = -49.308755760368
$ws.Range("N18").Value = -92.511912865895
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 15.384615384615
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = 8.474576271186
$ws.Range("I19").Value = 609
$ws.Range("J19").Value = 708
$ws.Range("K19").Value = -13.983050847457
$ws.Range("L19").Value = 58.181818181818
$ws.Range("M19").Value = 70.111731843575
$ws.Range("N19").Value = -43.715341959334
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 44
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 46.666666666666
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 229
$ws.Range("K20").Value = 8.733624454148
$ws.Range("L20").Value = 63.815789473684
$ws.Range("M20").Value = 46.470588235294
$ws.Range("N20").Value = -84.798534798534
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = 2.5
$ws.Range("F21").Value = 188
$ws.Range("G21").Value = 171
$ws.Range("H21").Value = 9.941520467836
$ws.Range("I21").Value = 1583
$ws.Range("J21").Value = 1612
$ws.Range("K21").Value = -1.799007444168
$ws.Range("L21").Value = 43.517679057117
$ws.Range("M21").Value = 24.743892828999
$ws.Range("N21").Value = -71.374321880651
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 68
$ws.Range("J22").Value = 46
$ws.Range("K22").Value = 47.826086956521
$ws.Range("L22").Value = 209.090909090909
$ws.Range("M22").Value = 161.538461538462
$ws.Range("C24").Value = 45
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 66.666666666666
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 137
$ws.Range("H24").Value = 3.649635036496
$ws.Range("I24").Value = 1449
$ws.Range("J24").Value = 1367
$ws.Range("K24").Value = 5.998536942209
$ws.Range("L24").Value = 46.068548387096
$ws.Range("M24").Value = 75.849514563106
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 88
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = 39.682539682539
$ws.Range("I25").Value = 714
$ws.Range("J25").Value = 649
$ws.Range("K25").Value = 10.015408320493
$ws.Range("L25").Value = 17.627677100494
$ws.Range("M25").Value = 3.179190751445
$ws.Range("D26").Value = 1
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("J26").Value = 39
$ws.Range("K26").Value = -12.820512820512
$ws.Range("L26").Value = -2.857142857142
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -75
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 30
$ws.Range("I27").Value = 119
$ws.Range("J27").Value = 77
$ws.Range("K27").Value = 54.545454545454
$ws.Range("L27").Value = 52.564102564102
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = 100
$ws.Range("L30").Value = 60
